$wb = $excel.ActiveWorkbook

# --- AR sheet ---
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = [double]"-8.196949479999377E-05"
$ws.Range("B3").Value = 0.7919696721424998
$ws.Range("B4").Value = 0.002073834418975113
$ws.Range("B5").Value = "[1.0, 0.32079848306623565, 0.3203464742749374, 0.30383664767653634, 0.2798568313288857, 0.4611029568298524, 0.28803404419654394, 0.30594939986523295, 0.28490971784891517, 0.2617786636679555, 0.24496516110730146, 0.23195750653055858, 0.2598396651269607, 0.25286424285428505, 0.26973344322086434, 0.2583382816011391, 0.20679644239414074, 0.23363571719323614, 0.20921430010915698, 0.2186498618795144]"

# --- SETAR sheet ---
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B3").Value = 1.060773480662983
$ws.Range("B4").Value = -0.01227695791452814
$ws.Range("B5").Value = 0.746199289076375
$ws.Range("B6").Value = 0.001246255001583923
$ws.Range("B7").Value = 0.01281216770763816
$ws.Range("B8").Value = 0.7464909570823748
$ws.Range("B9").Value = 0.001525397501077416
$ws.Range("B10").Value = "[1.0, 0.34846664012182754, 0.3247377283461004, 0.32978274763561793, 0.27778770719023343, 0.340320782054213, 0.2762465449480141, 0.27654745747550985, 0.2727665982663586, 0.24606579296650008, 0.2519179358383495, 0.21482364521630878, 0.23071921078354268, 0.22412884568464453, 0.2364145713556307, 0.2611590783395203, 0.20021590084537091, 0.21833260683478622, 0.2133731690708593, 0.20227042456278072]"

# --- GARCH sheet ---
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.0001931551601134447
$ws.Range("B3").Value = [double]"4.629545446446967E-05"
$ws.Range("B4").Value = 0.1000001338149829
$ws.Range("B5").Value = 0.8800000173411046
$ws.Range("B6").Value = "[1.0, -0.029494783234397005, -0.001969823308512687, -0.035374519020862624, -0.07400974511999941, 0.3054380810541253, -0.007307433594557866, 0.07733258989082864, 0.022415563674440912, -0.0018680692763889191, -0.02340574413524617, -0.019509466339076437, 0.02951933795402626, 0.03352462644888923, 0.08342770378050847, 0.027720743543569078, 0.002773267596367175, -0.009872712509838939, -0.01944839228819721, 0.022807578359101206]"

# --- TARCH sheet ---
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.0002847003374522831
$ws.Range("B3").Value = [double]"4.629535122323516E-05"
$ws.Range("B4").Value = 0.1000000631070222
$ws.Range("B5").Value = 0.854999999396984
$ws.Range("B6").Value = "[0.9999999999999998, -0.04690796695155648, -0.017069159133581577, -0.04883935806875373, -0.08655405172804705, 0.2953700681878833, -0.02111278575039498, 0.06658404820412295, 0.009475254108522613, -0.00967957905246629, -0.030363838633874936, -0.029302971300085905, 0.02127246399481523, 0.02444839814166219, 0.07922393999687773, 0.022075524532809546, -0.0037177421327436616, -0.015370411646061235, -0.025153865379800692, 0.019648753809765763]"
$ws.Range("B7").Value = 0.04999993404547573

# --- AR_TARCH sheet ---
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.0003038648332910364
$ws.Range("B3").Value = [double]"6.076152847966793E-06"
$ws.Range("B4").Value = 0.1303793987461055
$ws.Range("B5").Value = 0.9142725538503874
$ws.Range("B6").Value = "[0.9999999999999998, -0.02428935761627396, -0.004077367648715554, -0.0377347852749759, -0.04738805440505218, 0.195814960333635, -0.023038996144091124, 0.02844159635644579, -0.009820840067729205, -0.022027751982454505, -0.035656552275603784, -0.02482044457344369, -0.008377053577134483, -0.0028140405161985217, 0.04719432045927575, -0.0053434797688313995, -0.01817180717975269, -0.015500200642476963, -0.024335698102244535, -0.01611142641082071]"
$ws.Range("B7").Value = -0.08930390717180545
$ws.Range("B9").Value = 0.7683838256670201
